$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fifa_world_cup_2018_matches")

# Match 57 (row 58): status open -> completed, home goals 0 -> 2
$ws.Range("D58").Value = "completed"
$ws.Range("G58").Value = 2

# Match 58 (row 59): status open -> completed, home goals 0 -> 2, away goals 0 -> 1
$ws.Range("D59").Value = "completed"
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 1

# Match 59 (row 60): status open -> completed, home goals 0 -> 2
$ws.Range("D60").Value = "completed"
$ws.Range("G60").Value = 2

# Match 60 (row 61): status open -> completed, home goals 0 -> 1, away goals 0 -> 1
$ws.Range("D61").Value = "completed"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 1

# Match 61 (row 62): teams now known - Belgium vs France
$ws.Range("E62").Value = "Belgium"
$ws.Range("F62").Value = "France"

# Match 62 (row 63): teams now known - Croatia vs England
$ws.Range("E63").Value = "Croatia"
$ws.Range("F63").Value = "England"

# Update selection to reflect new active cell
$ws.Range("E64").Select()
